$d = $word.ActiveDocument

# --- Edit 1: "potrebno je da postoji" -> "nije" ---
$d.Content.Find.Execute("potrebno je da postoji", $true, $false, $false, $false, $false, $true, 1, $false, "nije", 2) | Out-Null

# --- Edit 2: "dugme" -> "odredjeno" ---
$d.Content.Find.Execute("dugme", $true, $false, $false, $false, $false, $true, 1, $false, "odredjeno", 2) | Out-Null

# --- Edit 3: "kojim" -> "cime" ---
$d.Content.Find.Execute("kojim", $true, $false, $false, $false, $false, $true, 1, $false, "cime", 2) | Out-Null

# --- Edit 4: "SSU document." -> "SSU dokument." ---
$d.Content.Find.Execute("SSU document.", $true, $false, $false, $false, $false, $true, 1, $false, "SSU dokument.", 2) | Out-Null
